$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRICES")   # "PRICES" is the active/tab-selected sheet

# New data rows to append after row 26 (row 27 = 2025-06-23, row 28 = 2025-07-08)
$rows = @(
    @{ Row = 27; Date = (Get-Date -Year 2025 -Month 6 -Day 23 -Hour 0 -Minute 0 -Second 0); Values = @(25.36, 4.5, 2.5, 12, 4.5, 6.2, 1.8, 1.599, 79.95, 11.95, 549.9, 49.95, 154.94999999999999, 64.95, 51.5) },
    @{ Row = 28; Date = (Get-Date -Year 2025 -Month 7 -Day 8  -Hour 0 -Minute 0 -Second 0); Values = @(26.17, 4.5, 2.5, 12, 4.5, 6.2, 1.8, 1.897, 119.95, 11.95, 599.95000000000005, 49.95, 154.94999999999999, 76.95, 51.21) }
)

foreach ($r in $rows) {
    $rowIndex = $r.Row

    # Column A: date, formatted/styled like the existing date column (same style as A3:A26)
    # Copy the format from the row above so the same shared style index is reused
    $ws.Cells.Item($rowIndex - 1, 1).Copy()
    $dateCell = $ws.Cells.Item($rowIndex, 1)
    $dateCell.PasteSpecial(-4122)  # xlPasteFormats
    $dateCell.Value = $r.Date

    # Columns B..P: numeric values
    for ($i = 0; $i -lt $r.Values.Count; $i++) {
        $col = $i + 2
        $ws.Cells.Item($rowIndex, $col).Value = $r.Values[$i]
    }
}

$excel.CutCopyMode = $false

# Update the active selection as recorded in the edited workbook
$ws.Range("N33").Select()
